$wb = $excel.ActiveWorkbook

# --- 1. "Fromage" sheet: add new row of cheese data (lot 002-204) ---
$wsFromage = $wb.Worksheets.Item("Fromage")
$wsFromage.Range("A8").Value = "002-204"
$wsFromage.Range("B8").Value = "Chèvre DC"
$wsFromage.Range("C8").Value = "C"
$wsFromage.Range("D8").Value = 3
$wsFromage.Range("E8").Value = 800
$wsFromage.Range("F8").Value = 4
[void]$wsFromage.Range("C13").Select()

# --- 2. "Lignes" sheet: fix a typo, add a "nom" column and fill in data ---
$wsLignes = $wb.Worksheets.Item("Lignes")

# fix trailing-space typo ("GF " -> "GF")
$wsLignes.Range("A4").Value = "GF"

# new "nom" header column, bold like the other header cells
$wsLignes.Range("A1").Value = "nom"
$wsLignes.Range("A1").Font.Bold = $true

# fill in the cadence / nombre d'employé / nombre d'opérateur columns
$wsLignes.Range("B2").Value = 5
$wsLignes.Range("C2").Value = 3
$wsLignes.Range("D2").Value = 1

$wsLignes.Range("B3").Value = 6
$wsLignes.Range("C3").Value = 2
$wsLignes.Range("D3").Value = 2

$wsLignes.Range("B4").Value = 1
$wsLignes.Range("C4").Value = 1
$wsLignes.Range("D4").Value = 2

$wsLignes.Range("B5").Value = 6
$wsLignes.Range("C5").Value = 3
$wsLignes.Range("D5").Value = 2

$wsLignes.Range("B6").Value = 3
$wsLignes.Range("C6").Value = 2
$wsLignes.Range("D6").Value = 1

# rename sheet (lowercase) and move it to the end, after "infos"
$wsLignes.Name = "lignes"
$wsInfos = $wb.Worksheets.Item("infos")
[void]$wsLignes.Move($null, $wsInfos)

# re-fetch fresh references after the move: the old variables can end up
# pointing at the wrong sheet since indices shift under Move()
$wsLignes = $wb.Worksheets.Item("lignes")
$wsInfos = $wb.Worksheets.Item("infos")

# --- 3. "infos" sheet: update the remembered selection ---
[void]$wsInfos.Range("B8:B9").Select()

# --- 4. "lignes" becomes the active sheet/tab, with its own selection ---
$wsLignes.Activate()
[void]$wsLignes.Range("E15").Select()
